# ---------------------------------------------------------------------------
# Applies the "Added functions for building constraints, as well as fixed
# some issues." commit to the workbook:
#   - processes sheet: new "is_cf_fix" column (C) + new "pv2" process row
#   - process_topology sheet: new "pv2 / sink / elc" topology row
#   - cf sheet: new pv2,s1 / pv2,s2 / pv2,s3 capacity-factor columns
#   - view/selection tweaks on nodes / processes / process_topology / cf
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsNodes   = $wb.Worksheets.Item(1)   # nodes
$wsProc    = $wb.Worksheets.Item(2)   # processes
$wsTopo    = $wb.Worksheets.Item(4)   # process_topology
$wsCf      = $wb.Worksheets.Item(6)   # cf

# ---------------------------------------------------------------------------
# 1) processes: insert the new "pv2" row (row 6) FIRST so the shared string
#    "pv2" is allocated before "pv2,s1"/"pv2,s2"/"pv2,s3"/"is_cf_fix".
# ---------------------------------------------------------------------------
$wsProc.Rows.Item(6).Insert()
$wsProc.Cells.Item(6, 1).Value = "pv2"
$wsProc.Cells.Item(6, 2).Value = 1
$wsProc.Cells.Item(6, 4).Value = 0
$wsProc.Cells.Item(6, 5).Value = 0
$wsProc.Cells.Item(6, 6).Value = 1
$wsProc.Cells.Item(6, 7).Value = 1
$wsProc.Cells.Item(6, 8).Value = 0
$wsProc.Cells.Item(6, 9).Value = 1
$wsProc.Cells.Item(6, 10).Value = 0
$wsProc.Cells.Item(6, 11).Value = 0
$wsProc.Cells.Item(6, 12).Value = 0

# ---------------------------------------------------------------------------
# 2) process_topology: insert the matching "pv2 / sink / elc" row (row 10),
#    reusing the "pv2" shared string allocated above.
# ---------------------------------------------------------------------------
$wsTopo.Rows.Item(10).Insert()
$wsTopo.Cells.Item(10, 1).Value = "pv2"
$wsTopo.Cells.Item(10, 2).Value = "sink"
$wsTopo.Cells.Item(10, 3).Value = "elc"
$wsTopo.Cells.Item(10, 4).Value = 1
$wsTopo.Cells.Item(10, 5).Value = 5
$wsTopo.Cells.Item(10, 6).Value = 0.5
$wsTopo.Cells.Item(10, 7).Value = 1
$wsTopo.Cells.Item(10, 8).Value = 1

# ---------------------------------------------------------------------------
# 3) cf: add pv2,s1 / pv2,s2 / pv2,s3 columns (E, F, G), mirroring pv1's
#    s1/s2/s3 columns (B, C, D). Column E used to be a blank, styled filler
#    column, so its leftover number format has to be cleared first.
# ---------------------------------------------------------------------------
$wsCf.Range("E2:G25").ClearFormats()

$wsCf.Cells.Item(1, 5).Value = "pv2,s1"
$wsCf.Cells.Item(1, 6).Value = "pv2,s2"
$wsCf.Cells.Item(1, 7).Value = "pv2,s3"

$cfValues = @(0, 0.4, 0.5, 0, 0.8, 1, 0.1, 0.6, 0.4, 0.6, 0.7, 0.1, 0.1, 0.8, 0.9, 0.2, 0.4, 0.6, 0.7, 0.7, 0.6, 0.7, 0.1, 0.6)
for ($i = 0; $i -lt $cfValues.Length; $i++) {
    $wsCf.Cells.Item($i + 2, 5).Value = $cfValues[$i]
}

$wsCf.Range("F2:F25").Formula = "=1*E2"
$wsCf.Range("G2:G25").Formula = "=1*E2"

# ---------------------------------------------------------------------------
# 4) processes: insert the "is_cf_fix" column (C) and its values, after the
#    "pv2" string has already claimed its shared-string slot.
# ---------------------------------------------------------------------------
$wsProc.Columns.Item(3).Insert()
$wsProc.Cells.Item(1, 3).Value = "is_cf_fix"
$wsProc.Cells.Item(2, 3).Value = 0
$wsProc.Cells.Item(3, 3).Value = 0
$wsProc.Cells.Item(4, 3).Value = 0
$wsProc.Cells.Item(5, 3).Value = 1
$wsProc.Cells.Item(6, 3).Value = 0
$wsProc.Cells.Item(7, 3).Value = 0
$wsProc.Cells.Item(8, 3).Value = 0

# ---------------------------------------------------------------------------
# 5) Selection / active-tab bookkeeping. The last sheet selected becomes the
#    active tab, so "processes" (the new active tab) is selected last.
# ---------------------------------------------------------------------------
$wsNodes.Range("C3").Select() | Out-Null
$wsTopo.Range("B15").Select() | Out-Null
$wsCf.Range("I8").Select() | Out-Null
$wsProc.Range("D10").Select() | Out-Null
